$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: updated timestamp text
$ws.Range("A1").Value = "Datos actualizados a 13 de Mayo de 2020 a las 23:35"

# Row 4: Estados Unidos - updated counts
$ws.Range("B4").Value = 1426320
$ws.Range("C4").Value = 17684
$ws.Range("D4").Value = 307296
$ws.Range("E4").Value = 1034085
$ws.Range("G4").Value = 1514
$ws.Range("H4").Value = 84939

# Rows 72-74: Camerun moves above Grecia/Azerbaiyan (reordered), with Camerun's stats refreshed
$ws.Range("A72").Value = "Camerun"
$ws.Range("B72").Value = 2800
$ws.Range("C72").Value = 111
$ws.Range("D72").Value = 1543
$ws.Range("E72").Value = 1121
$ws.Range("G72").Value = 11
$ws.Range("H72").Value = 136

$ws.Range("A73").Value = "Grecia"
$ws.Range("B73").Value = 2760
$ws.Range("C73").Value = 16
$ws.Range("D73").Value = 1374
$ws.Range("E73").Value = 1231
$ws.Range("F73").Value = 28
$ws.Range("G73").Value = 3
$ws.Range("H73").Value = 155

$ws.Range("A74").Value = "Azerbaiyan"
$ws.Range("B74").Value = 2758
$ws.Range("C74").Value = 65
$ws.Range("D74").Value = 1789
$ws.Range("E74").Value = 934
$ws.Range("F74").Value = 30
$ws.Range("G74").Value = 2
$ws.Range("H74").Value = 35

# Row 109: Guinea-Bisau - updated counts
$ws.Range("B109").Value = 836
$ws.Range("C109").Value = 16
$ws.Range("E109").Value = 807

# Row 148: Sudan del Sur - updated counts
$ws.Range("B148").Value = 203
$ws.Range("C148").Value = 9
$ws.Range("E148").Value = 201

# Rows 193-194: Nueva Caledonia moves above Belice (reordered, data swapped)
$ws.Range("A193").Value = "Nueva Caledonia"
$ws.Range("D193").Value = 18
$ws.Range("H193").Value = 0

$ws.Range("A194").Value = "Belice"
$ws.Range("D194").Value = 16
$ws.Range("H194").Value = 2
